$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$tcs = $s.ThemeColorScheme
$c = $tcs.Colors(5)
$c.RGB = 321
